# The data in rows 4-11 (Fecha, Calidad, Volumen, Precio minimo/maximo/promedio, Precio $/Kg)
# is cyclically shifted: each row's data moves up by 2 rows, wrapping rows 4-5 to the bottom (10-11).
# In other words, the new content of row R (for R in 4..11) equals the old content of row R+2
# (wrapping 10->4, 11->5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values actually change per the diff
$cols = @("D", "L", "M", "N", "O", "P", "S")

# Capture the original values for rows 4 through 11
$orig = @{}
for ($r = 4; $r -le 11; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $orig[$r] = $rowVals
}

# Mapping: new row -> old row supplying its data
$sourceRow = @{
    4  = 6
    5  = 7
    6  = 8
    7  = 9
    8  = 10
    9  = 11
    10 = 4
    11 = 5
}

foreach ($newRow in 4..11) {
    $oldRow = $sourceRow[$newRow]
    $vals = $orig[$oldRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value2 = $vals[$c]
    }
}
